# "version final sin errores" - update Version value and remove the
# Jurisdiction/Chile metadata row from the Metadata sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Version property value (row 3, column B): 0.4.0 -> 0.7.0
$ws.Range("B3").Value = "0.7.0"

# Remove the entire "Jurisdiction" / "Chile" row (row 11). Deleting the
# whole row shifts every row below it up by one, matching the new
# dimension (A1:B22 -> A1:B21) and re-indexing shared strings.
$ws.Rows.Item(11).Delete()
